# Weekly update to the "Hortaliza, Macroferia Regional de Talca - Alcachofa"
# price log: a new daily/weekly entry is inserted right after the existing
# entry for Madrigal/Primera (row 26), pushing every subsequent record down
# by one row. The sheet's used range grows from A1:R126 to A1:R127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 26 (shifts rows 26:126 down to 27:127).
$ws.Rows("26:26").Insert()

# Populate the newly inserted row with the new price-report record.
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 45145
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112013
$ws.Range("G26").Value = "Alcachofa"
$ws.Range("H26").Value = "Madrigal"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 12000
$ws.Range("N26").Value = "$/caja 40 unidades"
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 300
$ws.Range("Q26").Value = 40
$ws.Range("R26").Value = "Hortaliza"
